$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pr(>F) column (F) values re-rounded to match rebuilt PERMANOVA tables.
$updates = @{
    2  = 0.5914
    3  = 0.0014
    4  = 0.6889
    7  = 0.3631
    8  = 0.8555
    9  = 0.5987
    12 = 0.538
    13 = 0.9479
    14 = 0.2317
    17 = 0.783
    18 = 0.3419
    19 = 0.873
    22 = 0.154
    23 = 0.0969
    24 = 0.2408
    27 = 0.252
    28 = 0.6425
    29 = 0.6953
    32 = 0.8776
    33 = 0.0003
    34 = 0.9853
    37 = 0.3763
    38 = 0.0411
    39 = 0.4416
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
